
# IRGen: Add an SSA verifier pass to check SSA correctness (#81)
#
# Data-model changes represented in this edit:
#  - Stages sheet: rename the placeholder stage "TBD" -> "IRGen" (row 5)
#  - Errors sheet/table: append two new rows for the IRGen stage
#      Severity=Error, ID=1, Stage=IRGen  (no Description)
#      Severity=Error, ID=2, Stage=IRGen  (no Description)
#    with the same calculated ErrId / FullId formulas as the other rows.
#  - Table1 auto-expands from A1:F29 to A1:F31 to include the new rows.
#  - Selections move to follow the edit (Stages!A6, Errors!D31).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "TBD" placeholder stage to "IRGen" on the Stages sheet ---
$stages = $wb.Worksheets.Item("Stages")
$stages.Range("A5").Value = "IRGen"

# --- 2. Append the two new IRGen rows to the Errors table ---
$errors = $wb.Worksheets.Item("Errors")
$table = $errors.ListObjects.Item(1)

$row30 = $table.ListRows.Add()
$errors.Range("A30").Value = "Error"
$errors.Range("B30").Value = 1
$errors.Range("C30").Value = "IRGen"
$errors.Range("E30").Formula = '= (_xlfn.XLOOKUP($C30,Stages!$A:$A,Stages!$B:$B)+$B30)'
$errors.Range("F30").Formula = '= LEFT(A30,1)&E30'

$row31 = $table.ListRows.Add()
$errors.Range("A31").Value = "Error"
$errors.Range("B31").Value = 2
$errors.Range("C31").Value = "IRGen"
$errors.Range("E31").Formula = '= (_xlfn.XLOOKUP($C31,Stages!$A:$A,Stages!$B:$B)+$B31)'
$errors.Range("F31").Formula = '= LEFT(A31,1)&E31'

# --- 3. Recalculate so cached formula values are correct ---
$wb.Application.Calculate() | Out-Null

# --- 4. Match the resulting selections from the commit ---
$stages.Activate() | Out-Null
$stages.Range("A6").Select() | Out-Null

$errors.Activate() | Out-Null
$errors.Range("D31").Select() | Out-Null
